# feat: add 2022-Q3 data
#
# The workbook currently has two sheets:
#   1. "总计"      - summary of holdings per quarter
#   2. "2022-Q2"   - per-fund detail for the 2022-Q2 quarter
#
# This script:
#   1. Inserts a brand-new "2022-Q3" worksheet between them (so the
#      final sheet order is 总计, 2022-Q3, 2022-Q2) and fills it with the
#      2022-Q3 per-fund detail data (mirroring the layout of "2022-Q2").
#   2. Updates the "总计" summary sheet so that row 2 now reports the
#      2022-Q3 figures, and appends a new row 3 with the 2022-Q2 figures
#      that used to live in row 2.

$wb = $excel.ActiveWorkbook

$totalWs = $wb.Worksheets.Item("总计")
$q2Ws = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, positioned right before "2022-Q2"
# ---------------------------------------------------------------------
$q3Ws = $wb.Worksheets.Add($q2Ws)
$q3Ws.Name = "2022-Q3"

# Headers - same layout as the "2022-Q2" sheet
$q3Ws.Range("B1").Value = "基金代码"
$q3Ws.Range("C1").Value = "基金名称"
$q3Ws.Range("D1").Value = "基金规模"
$q3Ws.Range("E1").Value = "股票总仓位"
$q3Ws.Range("F1").Value = "仓位占比"
$q3Ws.Range("G1").Value = "持有市值(亿元)"
$q3Ws.Range("H1").Value = "仓位排名"

# Match the visual formatting (bold, centered, thin border) used on the
# equivalent header row of the "2022-Q2" sheet.
$q2Ws.Range("B1:H1").Copy()
$q3Ws.Range("B1:H1").PasteSpecial(-4122)

# Data row
$q3Ws.Range("A2").Value = 0

$q3Ws.Range("B2").NumberFormat = "@"
$q3Ws.Range("B2").Value = "014294"

$q3Ws.Range("C2").Value = "南方北交所精选两年定开混合"

$q3Ws.Range("D2").NumberFormat = "@"
$q3Ws.Range("D2").Value = "4.26"

$q3Ws.Range("E2").NumberFormat = "@"
$q3Ws.Range("E2").Value = "75.23"

$q3Ws.Range("F2").NumberFormat = "@"
$q3Ws.Range("F2").Value = "3.65"

$q3Ws.Range("G2").NumberFormat = "@"
$q3Ws.Range("G2").Value = "0.1555"

$q3Ws.Range("H2").Value = 6

# Match formatting on A2 too (bold/center/top/border like "2022-Q2"'s A2)
$q2Ws.Range("A2").Copy()
$q3Ws.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet
# ---------------------------------------------------------------------
# Row 2 now reflects 2022-Q3 instead of 2022-Q2
$totalWs.Range("B2").Value = "2022-Q3"
$totalWs.Range("D2").Value = 0.16

# Row 3 (new) holds the figures that used to be in row 2 for 2022-Q2
$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2022-Q2"
$totalWs.Range("C3").Value = 1
$totalWs.Range("D3").Value = 0.07

# Give A3 the same formatting as A2 (bold/center/top/border)
$totalWs.Range("A2").Copy()
$totalWs.Range("A3").PasteSpecial(-4122)
